$wb = $excel.ActiveWorkbook

# --- Sheet 1: EmitterLED ---
$ws1 = $wb.Worksheets.Item("EmitterLED")

# Update input values (dependent formula cells will recalculate automatically)
$ws1.Range("B3").Value = 0.19
$ws1.Range("B4").Value = 1.7

# Set the active cell selection to B7 on this sheet
$ws1.Activate()
$ws1.Range("B7").Select()

# --- Sheet 2: ReceiverTrans ---
$ws2 = $wb.Worksheets.Item("ReceiverTrans")

# Set the active cell selection to C10 on this sheet
$ws2.Activate()
$ws2.Range("C10").Select()
